# Add "Samples" and "Files" tab rows to the startup sheet,
# matching the other test workbooks in the suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab names first (keeps shared-string ordering consistent with the source file).
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# Row 3: SamplesTab
$ws.Range("B3").Value = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen IN  ["None"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("E3").Value = $ws.Range("E2").Value()

# Row 4: FilesTab
$ws.Range("B4").Value = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE tp.chemotherapy_regimen IN  ["None"]
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@
$ws.Range("C4").Value = $ws.Range("C2").Value()
$ws.Range("D4").Value = $ws.Range("D2").Value()
$ws.Range("E4").Value = $ws.Range("E2").Value()

# Match the wrap-text style already used for B2/C2 on the two new rows.
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

# Row heights match the wrapped-text content of each row.
$ws.Rows(2).RowHeight = 316.8
$ws.Rows(3).RowHeight = 345.6
$ws.Rows(4).RowHeight = 409.6

# Update the visible selection/scroll position like the saved workbook.
$ws.Range("C4").Select()
